# Fruta / hortaliza, semanal
#
# The weekly refresh re-sorted the daily price records (rows 2-19) into a
# new date order. The row contents themselves (Fecha / Volumen / Precio
# minimo / Precio maximo / Precio promedio ponderado / Precio $/Kg) did not
# change in value, they were simply redistributed across the rows, so we
# reproduce that by writing, for every row, the six values that belong to
# it in its new position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado) and P (Precio $/Kg)
# for rows 2 through 19, in order.
$rows = @(
    @{ Row = 2;  D = 44371; J = 34;  K = 5500; L = 6000; M = 5750; P = 359 },
    @{ Row = 3;  D = 44313; J = 34;  K = 6000; L = 6000; M = 6000; P = 375 },
    @{ Row = 4;  D = 44455; J = 52;  K = 5000; L = 6000; M = 5500; P = 344 },
    @{ Row = 5;  D = 44438; J = 34;  K = 5000; L = 6000; M = 5500; P = 344 },
    @{ Row = 6;  D = 44355; J = 25;  K = 6000; L = 6000; M = 6000; P = 375 },
    @{ Row = 7;  D = 44358; J = 52;  K = 6000; L = 6000; M = 6000; P = 375 },
    @{ Row = 8;  D = 44328; J = 160; K = 6000; L = 6000; M = 6000; P = 375 },
    @{ Row = 9;  D = 44442; J = 25;  K = 6000; L = 7000; M = 6480; P = 405 },
    @{ Row = 10; D = 44467; J = 52;  K = 5000; L = 6000; M = 5500; P = 344 },
    @{ Row = 11; D = 44308; J = 70;  K = 6000; L = 6000; M = 6000; P = 375 },
    @{ Row = 12; D = 44477; J = 25;  K = 6000; L = 6000; M = 6000; P = 375 },
    @{ Row = 13; D = 44306; J = 50;  K = 6000; L = 6000; M = 6000; P = 375 },
    @{ Row = 14; D = 44363; J = 160; K = 5500; L = 6000; M = 5750; P = 359 },
    @{ Row = 15; D = 44403; J = 43;  K = 6000; L = 6000; M = 6000; P = 375 },
    @{ Row = 16; D = 44474; J = 52;  K = 5000; L = 6000; M = 5500; P = 344 },
    @{ Row = 17; D = 44407; J = 45;  K = 5500; L = 6000; M = 5744; P = 359 },
    @{ Row = 18; D = 44341; J = 51;  K = 5500; L = 6000; M = 5755; P = 360 },
    @{ Row = 19; D = 44350; J = 25;  K = 6000; L = 6000; M = 6000; P = 375 }
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("J$r").Value = $entry.J
    $ws.Range("K$r").Value = $entry.K
    $ws.Range("L$r").Value = $entry.L
    $ws.Range("M$r").Value = $entry.M
    $ws.Range("P$r").Value = $entry.P
}
